# ---------------------------------------------------------------------------
# Insert a new "cost_rates" worksheet right before "process_mapping" and
# populate it with the parameter / unit / value table added by the commit
# "include cost rate in sheet".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# New sheet goes right before "process_mapping" -> matches target sheetId/order
$procMapping = $wb.Worksheets.Item("process_mapping")
$ws = $wb.Worksheets.Add($procMapping)
$ws.Name = "cost_rates"

# --- Header row -------------------------------------------------------------
$ws.Range("A1").Value = "parameter"
$ws.Range("B1").Value = "unit"
$ws.Range("C1").Value = "value"

# --- Data rows ---------------------------------------------------------------
$ws.Range("A2").Value = "variable_overhead_labor"
$ws.Range("B2").Value = "% of direct labor"
$ws.Range("C2").Value = 0.4

$ws.Range("A3").Value = "GSA_labour"
$ws.Range("B3").Value = "% of direct labor"
$ws.Range("C3").Value = 0.25

$ws.Range("A4").Value = "pack_profit"
$ws.Range("B4").Value = "% of investment"
$ws.Range("C4").Value = 0.05

$ws.Range("A5").Value = "launch_cost_labor"
$ws.Range("B5").Value = "% of direct labor"
$ws.Range("C5").Value = 0.1

$ws.Range("A6").Value = "launch_cost_material"
$ws.Range("B6").Value = "% of direct material cost"
$ws.Range("C6").Value = 0.05

$ws.Range("A7").Value = "working_capital"
$ws.Range("B7").Value = "% annual variable cost"
$ws.Range("C7").Value = 0.15

$ws.Range("A8").Value = "battery_warranty_costs"
$ws.Range("B8").Value = "% of pack cost"
$ws.Range("C8").Value = 0.056

$ws.Range("A9").Value = "variable_overhead_depreciation"
$ws.Range("B9").Value = "% of depreciation"
$ws.Range("C9").Value = 0.2

$ws.Range("A10").Value = "GSA_depreciation"
$ws.Range("B10").Value = "% of depreciation"
$ws.Range("C10").Value = 0.25

$ws.Range("A11").Value = "r_and_d"
$ws.Range("B11").Value = "% of depreciation"
$ws.Range("C11").Value = 0.4

$ws.Range("A12").Value = "lifetime_capital_equipment"
$ws.Range("B12").Value = "years"
$ws.Range("C12").Value = 10

# --- Formatting (Consolas / black / vertically centred) ---------------------
# Applied in the same grouping the original author used: first just column A
# for the top block, then full A:E for a couple of rows, then A:C for the
# bottom block.
$fmt = $ws.Range("A3:A6")
$fmt.Font.Color = 0
$fmt.Font.Name = "Consolas"
$fmt.VerticalAlignment = -4108

$fmt2 = $ws.Range("A7:E8")
$fmt2.Font.Color = 0
$fmt2.Font.Name = "Consolas"
$fmt2.VerticalAlignment = -4108

$fmt3 = $ws.Range("A9:C12")
$fmt3.Font.Color = 0
$fmt3.Font.Name = "Consolas"
$fmt3.VerticalAlignment = -4108

# Column widths to match the new sheet layout
$ws.Columns("A").ColumnWidth = 39.28515625
$ws.Columns("B").ColumnWidth = 26.7109375

# --- Selection bookkeeping ----------------------------------------------------
# Leave the "default_manufacturing_rates" sheet's own selection where the user
# had left it, then make the new sheet the active tab.
$dmr = $wb.Worksheets.Item("default_manufacturing_rates")
$dmr.Activate()
$dmr.Range("E11").Select()

$ws.Activate()
$ws.Range("G21").Select()

Write-Output "cost_rates sheet inserted"
